$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Journal": fill in a new journal entry on row 7
#   (work performed on PotenotTaskServece: tests for PotenotTaskCalculatorImpl)
# ---------------------------------------------------------------------------
$wsJournal = $wb.Worksheets.Item("Journal")

$wsJournal.Range("A7").Value = "Разработка тестов для PotenotService"
$wsJournal.Range("B7").Value = 45544
$wsJournal.Range("C7").Value = 0.0833333333333333
$wsJournal.Range("D7").Value = "PotenotTaskServece"

# ---------------------------------------------------------------------------
# Sheet "PotenotTask test data": fill in computed columns N:U for target2..4
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("PotenotTask test data")

# target2 (row 3)
$wsData.Range("N3").Value = 0.791317738
$wsData.Range("O3").Value = 45.2021
$wsData.Range("P3").Value = 1.579809013
$wsData.Range("Q3").Value = 90.3059
$wsData.Range("R3").Value = 625.62443
$wsData.Range("S3").Value = -1176.961644
$wsData.Range("T3").Value = 152.0024
$wsData.Range("U3").Value = 2.653017829

# target3 (row 4)
$wsData.Range("N4").Value = 0.127026032
$wsData.Range("O4").Value = 7.1641
$wsData.Range("P4").Value = 0.151140665
$wsData.Range("Q4").Value = 8.3935
$wsData.Range("T4").Value = 317.1807
$wsData.Range("U4").Value = 5.537963654

# target4 (row 5)
$wsData.Range("N5").Value = 0.163236954
$wsData.Range("O5").Value = 20.4843
$wsData.Range("P5").Value = 0.754787027
$wsData.Range("Q5").Value = 43.1446
$wsData.Range("T5").Value = 163.0559
$wsData.Range("U5").Value = 2.846627162

# ---------------------------------------------------------------------------
# View state: zoom in to 120% on both sheets and move the selection down,
# matching where the author continued working.
# ---------------------------------------------------------------------------
$wsData.Activate()
$excel.ActiveWindow.Zoom = 120
$wsData.Range("U3").Select()

$wsJournal.Activate()
$excel.ActiveWindow.Zoom = 120
$wsJournal.Range("A29").Select()
